# Scheduled runner update: refresh computed leve-profit figures (Cerberus_Profits)
# across all job sheets. Each block below updates the price/profit columns
# (H..N) for the specific rows whose underlying market data changed.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 114
$ws.Range("I9").Value = 108.25
$ws.Range("J9").Value = 125.5
$ws.Range("K9").Value = 108.25
$ws.Range("L9").Value = 125.5
$ws.Range("M9").Value = 60.75
$ws.Range("N9").Value = -463.5
# Row 28
$ws.Range("H28").Value = 338.75
$ws.Range("I28").Value = 356.25
$ws.Range("J28").Value = 268.75
$ws.Range("K28").Value = 356.25
$ws.Range("L28").Value = 268.75
$ws.Range("M28").Value = 128.75
$ws.Range("N28").Value = -1238.75
# Row 33
$ws.Range("H33").Value = 1292.8889
$ws.Range("I33").Value = 204.25
$ws.Range("K33").Value = 204.25
$ws.Range("M33").Value = 24.75
# Row 43
$ws.Range("H43").Value = 2299.6
$ws.Range("I43").Value = 2199.5
$ws.Range("K43").Value = 2199.5
$ws.Range("M43").Value = -2130.5
# Row 58
$ws.Range("H58").Value = 2274.0908
$ws.Range("I58").Value = 1001.5
$ws.Range("K58").Value = 3004.5
$ws.Range("M58").Value = -2854.5
# Row 95
$ws.Range("H95").Value = 33395.6
$ws.Range("J95").Value = 33395.6
$ws.Range("L95").Value = 33395.6
$ws.Range("N95").Value = -38887.6
# Row 96
$ws.Range("H96").Value = 762
$ws.Range("I96").Value = 725.75
$ws.Range("J96").Value = 798.25
$ws.Range("K96").Value = 2177.25
$ws.Range("L96").Value = 2394.75
$ws.Range("M96").Value = -804.25
$ws.Range("N96").Value = -5140.75
# Row 106
$ws.Range("H106").Value = 4800591.5
$ws.Range("I106").Value = 5696965
$ws.Range("K106").Value = 5696965
$ws.Range("M106").Value = -5696334
# Row 113
$ws.Range("H113").Value = 7113.615
$ws.Range("I113").Value = 6838.5557
$ws.Range("J113").Value = 7259.2354
$ws.Range("K113").Value = 6838.5557
$ws.Range("L113").Value = 7259.2354
$ws.Range("M113").Value = -3584.5557
$ws.Range("N113").Value = -13767.2354
# Row 116
$ws.Range("H116").Value = 16340.5
$ws.Range("I116").Value = 18573.834
$ws.Range("J116").Value = 14107.167
$ws.Range("K116").Value = 18573.834
$ws.Range("L116").Value = 14107.167
$ws.Range("M116").Value = -15131.834
$ws.Range("N116").Value = -20991.167
# Row 125
$ws.Range("H125").Value = 3037.625
$ws.Range("I125").Value = 1948.6
$ws.Range("J125").Value = 4852.6665
$ws.Range("K125").Value = 17537.4
$ws.Range("L125").Value = 43673.9985
$ws.Range("M125").Value = -15077.4
$ws.Range("N125").Value = -48593.9985
# Row 127
$ws.Range("H127").Value = 1137.6111
$ws.Range("I127").Value = 1038.6
$ws.Range("K127").Value = 3115.8
$ws.Range("M127").Value = 1844.2
# Row 132
$ws.Range("H132").Value = 4295.8975
$ws.Range("I132").Value = 4581.9033
$ws.Range("J132").Value = 3187.625
$ws.Range("K132").Value = 13745.7099
$ws.Range("L132").Value = 9562.875
$ws.Range("M132").Value = -11215.7099
$ws.Range("N132").Value = -14622.875
# Row 134
$ws.Range("H134").Value = 42996.785
$ws.Range("J134").Value = 42996.785
$ws.Range("L134").Value = 42996.785
$ws.Range("N134").Value = -53136.785
# Row 135
$ws.Range("H135").Value = 1263.35
$ws.Range("I135").Value = 1181.5555
$ws.Range("K135").Value = 10633.9995
$ws.Range("M135").Value = -8098.9995
# Row 137
$ws.Range("H137").Value = 402304.28
$ws.Range("I137").Value = 626650.9
$ws.Range("J137").Value = 3465.889
$ws.Range("K137").Value = 1879952.7
$ws.Range("L137").Value = 10397.667
$ws.Range("M137").Value = -1877402.7
$ws.Range("N137").Value = -15497.667

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 30
$ws.Range("H30").Value = 8500
$ws.Range("I30").Value = 10000
$ws.Range("K30").Value = 10000
$ws.Range("M30").Value = -9850
# Row 32
$ws.Range("H32").Value = 2633.59
$ws.Range("I32").Value = 2174.6792
$ws.Range("K32").Value = 2174.6792
$ws.Range("M32").Value = -1887.6792
# Row 74
$ws.Range("H74").Value = 1832.5938
$ws.Range("I74").Value = 1097
$ws.Range("J74").Value = 2666.2666
$ws.Range("K74").Value = 1097
$ws.Range("L74").Value = 2666.2666
$ws.Range("M74").Value = -223
$ws.Range("N74").Value = -4414.2666
# Row 77
$ws.Range("H77").Value = 1832.5938
$ws.Range("I77").Value = 1097
$ws.Range("J77").Value = 2666.2666
$ws.Range("K77").Value = 5485
$ws.Range("L77").Value = 13331.333
$ws.Range("M77").Value = -1117
$ws.Range("N77").Value = -22067.333
# Row 95
$ws.Range("H95").Value = 29666.334
$ws.Range("J95").Value = 29666.334
$ws.Range("L95").Value = 29666.334
$ws.Range("N95").Value = -35158.334
# Row 122
$ws.Range("H122").Value = 2106.125
$ws.Range("I122").Value = 2210.3635
$ws.Range("J122").Value = 1876.8
$ws.Range("K122").Value = 6631.0905
$ws.Range("L122").Value = 5630.4
$ws.Range("M122").Value = -4181.0905
$ws.Range("N122").Value = -10530.4
# Row 132
$ws.Range("H132").Value = 2941.3572
$ws.Range("I132").Value = 2498.4614
$ws.Range("J132").Value = 8699
$ws.Range("K132").Value = 7495.3842
$ws.Range("L132").Value = 26097
$ws.Range("M132").Value = -4965.3842
$ws.Range("N132").Value = -31157
# Row 138
$ws.Range("H138").Value = 99999
$ws.Range("I138").Value = 100000
$ws.Range("K138").Value = 100000
$ws.Range("M138").Value = -94860

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 36465.5
$ws.Range("I26").Value = 31758.6
$ws.Range("K26").Value = 31758.6
$ws.Range("M26").Value = -31466.6
# Row 105
$ws.Range("H105").Value = 1529.3684
$ws.Range("I105").Value = 1717.9286
$ws.Range("K105").Value = 1717.9286
$ws.Range("M105").Value = 29.07140000000004
# Row 134
$ws.Range("H134").Value = 9996.24
$ws.Range("I134").Value = 8732.736999999999
$ws.Range("J134").Value = 13997.333
$ws.Range("K134").Value = 26198.211
$ws.Range("L134").Value = 41991.999
$ws.Range("M134").Value = -23663.211
$ws.Range("N134").Value = -47061.999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 11933.871
$ws.Range("J4").Value = 11933.871
$ws.Range("L4").Value = 11933.871
$ws.Range("N4").Value = -12157.871
# Row 105
$ws.Range("H105").Value = 1248.9166
$ws.Range("I105").Value = 1099.3684
$ws.Range("J105").Value = 1817.2
$ws.Range("K105").Value = 1099.3684
$ws.Range("L105").Value = 1817.2
$ws.Range("M105").Value = 647.6315999999999
$ws.Range("N105").Value = -5311.2
# Row 107
$ws.Range("H107").Value = 3603.1667
$ws.Range("I107").Value = 3910.25
$ws.Range("K107").Value = 3910.25
$ws.Range("M107").Value = -1990.25
# Row 124
$ws.Range("H124").Value = 88000
$ws.Range("J124").Value = 88000
$ws.Range("L124").Value = 88000
$ws.Range("N124").Value = -92910

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 18
$ws.Range("H18").Value = 264
$ws.Range("I18").Value = 99.666664
$ws.Range("J18").Value = 428.33334
$ws.Range("K18").Value = 298.999992
$ws.Range("L18").Value = 1285.00002
$ws.Range("M18").Value = -129.999992
$ws.Range("N18").Value = -1623.00002
# Row 99
$ws.Range("H99").Value = 3824.8333
$ws.Range("I99").Value = 589.8
$ws.Range("K99").Value = 1769.4
$ws.Range("M99").Value = 476.6000000000001
# Row 113
$ws.Range("H113").Value = 1624.4445
$ws.Range("I113").Value = 1196.1666
$ws.Range("J113").Value = 2481
$ws.Range("K113").Value = 3588.4998
$ws.Range("L113").Value = 7443
$ws.Range("M113").Value = -1418.4998
$ws.Range("N113").Value = -11783
# Row 114
$ws.Range("H114").Value = 1186
$ws.Range("J114").Value = 1299.6666
$ws.Range("L114").Value = 3898.9998
$ws.Range("N114").Value = -10406.9998
# Row 129
$ws.Range("H129").Value = 27781204
$ws.Range("J129").Value = 33337302
$ws.Range("L129").Value = 100011906
$ws.Range("N129").Value = -100021906

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 750.2632
$ws.Range("J97").Value = 1371.1111
$ws.Range("L97").Value = 1371.1111
$ws.Range("N97").Value = -2363.1111
# Row 102
$ws.Range("H102").Value = 4467.25
$ws.Range("I102").Value = 4491.1055
$ws.Range("K102").Value = 4491.1055
$ws.Range("M102").Value = -2869.1055
# Row 122
$ws.Range("H122").Value = 2591.25
$ws.Range("I122").Value = 1455.3334
$ws.Range("K122").Value = 4366.0002
$ws.Range("M122").Value = -1916.0002
# Row 123
$ws.Range("H123").Value = 87500
$ws.Range("J123").Value = 87500
$ws.Range("L123").Value = 87500
$ws.Range("N123").Value = -92400
# Row 132
$ws.Range("H132").Value = 4136.478
$ws.Range("I132").Value = 4072.35
$ws.Range("J132").Value = 4564
$ws.Range("K132").Value = 12217.05
$ws.Range("L132").Value = 13692
$ws.Range("M132").Value = -9687.049999999999
$ws.Range("N132").Value = -18752
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3086.9167
$ws.Range("I7").Value = 3117.5454
$ws.Range("K7").Value = 3117.5454
$ws.Range("M7").Value = -3005.5454
# Row 16
$ws.Range("H16").Value = 1754.5
$ws.Range("I16").Value = 1754.5
$ws.Range("K16").Value = 1754.5
$ws.Range("M16").Value = -1584.5
# Row 126
$ws.Range("H126").Value = 3086.9167
$ws.Range("I126").Value = 3117.5454
$ws.Range("K126").Value = 9352.636200000001
$ws.Range("M126").Value = -6882.636200000001
